# BUG FIX 1.) BTP SINGLE EXPORT  2.) DSR FROM VALIDATION
#
# The "Resources >>" block on the BTPReport sheet used to list three
# resource names (Suresh Kannan / Prem Anandakrishnan / Saranya Rajendran)
# in the merged cells D2:E2, F2:G2 and D3:E3. Those stray/hard-coded
# names were removed so the template ships with a clean, blank
# "Resources >>" area that gets populated at export time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("BTPReport")

# Clear the three hard-coded resource names (merged ranges keep their
# merge/formatting, only the text content is removed).
$ws.Range("D2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("D3").Value = ""

# Reflect the selection that was active when the sheet was last saved.
$ws.Range("B10:C10").Select()
